$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Hoja1")

# Update the date cell (A1)
$ws.Range("A1").Value = 45311

# Update price values in column D (bug fix for exceeded request in google drive pricing sync)
$updates = @{
    "D24" = 617.965
    "D25" = 638.478
    "D26" = 656.423
    "D27" = 700.783
    "D28" = 853.096
    "D29" = 1000.023
    "D30" = 1143.613
    "D31" = 1235.924
    "D36" = 791.045
    "D37" = 868.736
    "D38" = 920.0170000000001
    "D39" = 975.664
    "D40" = 1089.769
    "D41" = 1415.417
    "D42" = 1430.801
    "D43" = 1466.7
    "D44" = 1871.839
    "D45" = 2197.487
    "D46" = 2692.371
    "D47" = 3082.117
    "D52" = 843.61
    "D53" = 1015.406
    "D54" = 1115.411
    "D55" = 1215.414
    "D56" = 1320.543
    "D57" = 1466.7
    "D58" = 1576.961
    "D59" = 1710.297
    "D60" = 2527.749
    "D61" = 2735.962
    "D62" = 3043.659
    "D63" = 3651.367
    "D68" = 1013.358
    "D69" = 1165.152
    "D70" = 1215.414
    "D71" = 1320.543
    "D72" = 1503.37
    "D73" = 1625.677
    "D74" = 1782.092
    "D75" = 1980.042
    "D76" = 2676.983
    "D77" = 3138.533
    "D78" = 3674.439
    "D79" = 3961.624
    "D85" = 1416.697
    "D86" = 1646.193
    "D87" = 1836.453
    "D88" = 2023.122
    "D89" = 2425.697
    "D90" = 2539.032
    "D91" = 2994.939
    "D92" = 3359.049
    "D93" = 3974.448
    "D94" = 4484.714
    "D95" = 5479.608
    "D101" = 2143.638
    "D102" = 2425.697
    "D103" = 2738.521
    "D104" = 3246.228
    "D105" = 3282.125
    "D106" = 4146.244
    "D107" = 4666.772
    "D108" = 5269.353
    "D109" = 6102.7
    "D115" = 3075.712
    "D116" = 3400.078
    "D117" = 3948.807
    "D118" = 4100.091
    "D119" = 4615.485
    "D120" = 5379.607
    "D121" = 6282.192
    "D122" = 7102.722
    "D123" = 7718.12
}

foreach ($cell in $updates.Keys) {
    $ws.Range($cell).Value = $updates[$cell]
}
